$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price and Volume columns keep their literal text formatting
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "29.181.63"
$ws.Range("E2").Value = "  -0.93%  "
$ws.Range("D3").Value = "1.862.02"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "0.7076"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "240.91"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -0.80%  "
$ws.Range("D9").Value = "0.07647"
$ws.Range("E9").Value = "  -2.37%  "
$ws.Range("D10").Value = "24.67"
$ws.Range("E10").Value = "  -1.95%  "
$ws.Range("D11").Value = "0.08326"
$ws.Range("E11").Value = "  +0.99%  "
$ws.Range("D12").Value = "1.843.00"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").Value = "5.179"
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "0.7096"
$ws.Range("E14").Value = "  -2.46%  "
$ws.Range("D15").Value = "91.21"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("D16").Value = "29.209.78"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("D17").Value = "5.917"
$ws.Range("E17").Value = "  -0.25%  "
$ws.Range("D18").Value = "242.64"
$ws.Range("D19").Value = "0.000007814"
$ws.Range("E19").Value = "  -0.76%  "
$ws.Range("D20").Value = "2.114.81"
$ws.Range("E20").Value = "  +0.09%  "
$ws.Range("D21").Value = "13.07"
$ws.Range("E21").Value = "  -1.65%  "
$ws.Range("D22").Value = "0.9999"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "7.855"
$ws.Range("E23").Value = "  -1.46%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").Value = "0.1592"
$ws.Range("E25").Value = "  +0.86%  "
$ws.Range("D26").Value = "163.41"
$ws.Range("E26").Value = "  -0.20%  "
$ws.Range("D27").Value = "8.934"
$ws.Range("E27").Value = "  -0.89%  "
$ws.Range("D28").Value = "18.46"
$ws.Range("E28").Value = "  +0.80%  "
$ws.Range("E29").Value = "  -0.24%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "1.313"
$ws.Range("E30").Value = "  -3.79%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "4.406"
$ws.Range("E31").Value = "  +0.54%  "
$ws.Range("D32").Value = "4.221"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("D33").Value = "0.05130"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("D34").Value = "0.7972"
$ws.Range("E34").Value = "  +10.17%  "
$ws.Range("D35").Value = "1.908"
$ws.Range("E35").Value = "  -1.36%  "
$ws.Range("D36").Value = "1.161"
$ws.Range("E36").Value = "  -3.28%  "
$ws.Range("D37").Value = "2.681"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("D38").Value = "0.01840"
$ws.Range("E38").Value = "  -1.34%  "
$ws.Range("D39").Value = "2.697"
$ws.Range("E39").Value = "  -1.10%  "
$ws.Range("D40").Value = "1.166.71"
$ws.Range("E40").Value = "  -6.05%  "
$ws.Range("D41").Value = "6.178"
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "0.8884"
$ws.Range("E42").Value = "  -2.12%  "
$ws.Range("D43").Value = "72.82"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "102.14"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").Value = "2.013.16"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").Value = "0.5187"
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("D48").Value = "1.768"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "9.306"
$ws.Range("E51").Value = "  -0.02%  "
